# Update "want-to-go" (想去人数) counts on sheet "展览" and "全部类型"
# to reflect the latest generated data (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Map of row number -> new value for sheet "展览"
$updatesExpo = @{
    3  = 7871
    8  = 604
    13 = 3144
    21 = 247
    23 = 312
    28 = 19
    32 = 516
    37 = 95
}

$wsExpo = $wb.Worksheets.Item("展览")
foreach ($row in $updatesExpo.Keys) {
    $wsExpo.Range("F$row").Value = $updatesExpo[$row]
}

# Map of row number -> new value for sheet "全部类型"
$updatesAll = @{
    5  = 7871
    10 = 604
    16 = 3144
    26 = 247
    28 = 312
    33 = 19
    37 = 516
    42 = 95
}

$wsAll = $wb.Worksheets.Item("全部类型")
foreach ($row in $updatesAll.Keys) {
    $wsAll.Range("F$row").Value = $updatesAll[$row]
}
